$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new daily rows after the existing data (which ends at row 74).
$ws.Cells.Item(75, 1).Value = 46037
$ws.Cells.Item(75, 2).Value = 731
$ws.Cells.Item(75, 3).Value = 695
$ws.Cells.Item(75, 4).Value = 36

$ws.Cells.Item(76, 1).Value = 46038
$ws.Cells.Item(76, 2).Value = 679
$ws.Cells.Item(76, 3).Value = 669
$ws.Cells.Item(76, 4).Value = 10

# The new date cells get a (new) date number format, distinct from the
# mm-dd-yy format used by the rest of column A.
$ws.Range("A75:A76").NumberFormat = "d-mmm-yy"

# Match the author's final selection/view state.
$ws.Range("A75:D76").Select()
